# Management_information_KAL_2016.xlsx - align values with the commit:
# "Updated multiple management files to have consistent and english values."
#
# Concretely (sheet "BON_2016"):
#   - E44:E46 ("Notice" for the low-N/no-fungicide treatment 1st/2nd/3rd N
#     applications) had a stray literal 0 value; cleared to blank.
#   - D56:D58 and D68:D70 ("Amount" for the 1st N application of the
#     high-N treatments) stored the amount as text ("50 kg/ha", "55 kg/ha",
#     "60 kg/ha"); changed to plain numeric amounts (unit now implicit).
#   - E56:E58 and E68:E70 ("Notice" for those same rows) held the German
#     fertilizer name "N-Düngung KAS"; replaced with the English/standard
#     abbreviation "CAN" (Calcium Ammonium Nitrate).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BON_2016")

# --- Clear the stray 0s in E44:E46 (keep existing formatting) ---
$ws.Range("E44").ClearContents()
$ws.Range("E45").ClearContents()
$ws.Range("E46").ClearContents()

# --- D56 / D68: text "50 kg/ha" -> number 50, keeping the cell's
#     original (quote-prefixed) number format/style ---
foreach ($addr in @("D56", "D68")) {
    $cell = $ws.Range($addr)
    $cell.Copy() | Out-Null
    $ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cell.Value = 50
    $ws.Range("ZZ1").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null              # xlPasteFormats
}
$ws.Range("ZZ1").Clear()

# --- D57 / D69: text "55 kg/ha" -> number 55 ---
$ws.Range("D57").Value = 55
$ws.Range("D69").Value = 55

# --- D58 / D70: text "60 kg/ha" -> number 60 ---
$ws.Range("D58").Value = 60
$ws.Range("D70").Value = 60

# --- E56:E58 / E68:E70: "N-Düngung KAS" -> "CAN" ---
foreach ($addr in @("E56", "E57", "E58", "E68", "E69", "E70")) {
    $ws.Range($addr).Value = "CAN"
}
